$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: all zero for rows 2-9
$ws.Range("B2:B9").Value = 0

# Column C values
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = -0.7222377876285779
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 0.6810405984895163
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = -0.6773077325108112
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = -0.6485656446443729

# Column D values
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 0.7528533129576779
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = -0.7387131370077278
$ws.Range("D8").Value = 0
$ws.Range("D9").Value = 0
